$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Right" counts in the Marking and Total rows,
# and the corresponding Correct/Total text in the Total row.
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 40
$ws.Range("E12").Value = "40/140"
